$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "RandomPassword: DyFbkobqRandomEmail: rMETeA@gmailAccountCreated",
    "RandomPassword: DyFbkobqRandomEmail: rMETeA@gmailAccountDeleted",
    "RandomPassword: hhoSOq2oRandomEmail: zymyvx@gmailAccountCreated",
    "RandomPassword: Rl5D7YakRandomEmail: YkvChS@gmailAccountCreated",
    "RandomPassword: AQ8BNTt2RandomEmail: wLHFlN@gmailAccountCreated",
    "RandomPassword: AQ8BNTt2RandomEmail: wLHFlN@gmailAccountDeleted",
    "RandomPassword: lhWAs14ERandomEmail: XHKQrv@gmailAccountCreated"
)

$startRow = 71
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
